$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "{'shimenet', 'naur', 'andamhie', 'anda', 'chika', 'eklabool'}"
$ws.Range("D3").Value = "{'shimenet', 'naur', 'andamhie', 'anda', 'chika', 'eklabool', 'λ'}"
$ws.Range("D4").Value = "{'shimenet', 'naur', 'andamhie', 'anda', 'chika', 'eklabool', 'λ'}"
$ws.Range("D6").Value = "{'andamhie', 'shimenet', 'anda', 'chika', 'eklabool'}"
$ws.Range("D7").Value = "{'andamhie', 'eklabool', 'chika', 'anda'}"
$ws.Range("D10").Value = "{'chika', 'andamhie', 'anda', 'eklabool', 'λ'}"
$ws.Range("D11").Value = "{'andamhie', 'eklabool', 'anda', 'chika'}"
$ws.Range("D12").Value = "{'λ', ','}"
$ws.Range("D17").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D18").Value = "{'andamhie_literal', 'chika_literal', 'korik', 'eme', 'anda_literal'}"
$ws.Range("D19").Value = "{'eme', 'korik'}"
$ws.Range("D24").Value = "{'λ', ','}"
$ws.Range("D26").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D28").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '{', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D29").Value = "{'λ', ','}"
$ws.Range("D30").Value = "{'shimenet', 'andamhie', 'anda', 'chika', 'eklabool', 'λ'}"
$ws.Range("D31").Value = "{'forda', 'id', 'pak', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'andamhie', 'anda', 'keri', 'versa', 'serve', 'adele'}"
$ws.Range("D32").Value = "{'andamhie', 'naur', 'anda', 'chika', 'eklabool', 'λ'}"
$ws.Range("D33").Value = "{'naur', 'andamhie', 'anda', 'chika', 'eklabool'}"
$ws.Range("D34").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', 'λ', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D35").Value = "{'λ', ','}"
$ws.Range("D36").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D38").Value = "{'id', 'serve', 'pak', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'andamhie', 'anda', 'keri', 'versa', 'forda', 'adele'}"
$ws.Range("D39").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D40").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D42").Value = "{'*', '+', '>=', '/', '==', '%', '<=', 'λ', '//', '**', '>', '!=', '&&', '-', '<', '||'}"
$ws.Range("D43").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D44").Value = "{'-', '!', 'λ'}"
$ws.Range("D45").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', 'chika_literal', '--', '++', 'anda_literal'}"
$ws.Range("D47").Value = "{'(', '[', '++', '--', 'λ'}"
$ws.Range("D49").Value = "{'andamhie_literal', 'chika_literal', 'korik', 'eme', 'anda_literal'}"
$ws.Range("D50").Value = "{'*', '+', '>=', '/', '==', '%', '<=', '//', '**', '>', '!=', '&&', '-', '<', '||'}"
$ws.Range("D51").Value = "{'id', 'serve', 'pak', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'andamhie', 'anda', 'keri', 'versa', 'forda', 'adele'}"
$ws.Range("D52").Value = "{'id', 'serve', 'pak', 'gogogo', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'amaccana', 'andamhie', 'anda', 'keri', 'versa', 'forda', 'adele'}"
$ws.Range("D53").Value = "{'forda', 'id', 'pak', 'gogogo', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'amaccana', 'andamhie', 'anda', 'keri', 'versa', 'serve', 'adele'}"
$ws.Range("D54").Value = "{'(', '/=', '+=', '*=', '**=', '[', '=', '//=', '-=', '%='}"
$ws.Range("D55").Value = "{'*=', '=', '**=', '+=', '/=', '//=', '-=', '%='}"
$ws.Range("D56").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '{', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D57").Value = "{'andamhie', 'id', 'anda', 'chika', 'eklabool'}"
$ws.Range("D58").Value = "{'chika', 'andamhie', 'anda', 'eklabool', 'λ'}"
$ws.Range("D59").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D63").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D65").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D66").Value = "{'forda', 'id', 'pak', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'andamhie', 'anda', 'keri', 'versa', 'serve', 'adele'}"
$ws.Range("D68").Value = "{'forda', 'id', 'pak', 'gogogo', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'amaccana', 'andamhie', 'anda', 'keri', 'versa', 'serve', 'adele'}"
$ws.Range("D69").Value = "{'λ', 'ganern'}"
$ws.Range("D70").Value = "{'λ', 'ganern'}"
$ws.Range("D73").Value = "{'chika', 'andamhie', 'anda', 'eklabool', 'λ'}"
$ws.Range("D74").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D75").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D77").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', 'λ', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D78").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D79").Value = "{'forda', 'id', 'pak', 'gogogo', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'amaccana', 'andamhie', 'anda', 'keri', 'versa', 'serve', 'adele'}"
$ws.Range("D81").Value = "{'lang', '('}"
$ws.Range("D82").Value = "{'forda', 'id', 'pak', 'gogogo', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'amaccana', 'andamhie', 'anda', 'keri', 'versa', 'serve', 'adele'}"
$ws.Range("D84").Value = "{'betsung', 'λ'}"
$ws.Range("D85").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
$ws.Range("D86").Value = "{'forda', 'id', 'pak', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'andamhie', 'anda', 'keri', 'versa', 'serve', 'adele'}"
$ws.Range("D88").Value = "{'betsung', 'λ'}"
$ws.Range("D89").Value = "{'forda', 'id', 'pak', 'gogogo', '--', 'chika', 'push', 'eklabool', '++', 'λ', 'naur', 'adelete', 'amaccana', 'andamhie', 'anda', 'keri', 'versa', 'serve', 'adele'}"
$ws.Range("D92").Value = "{'amaccana', 'gogogo', 'λ'}"
$ws.Range("D93").Value = "{'push', 'λ'}"
$ws.Range("D94").Value = "{'--', 'id', '++'}"
$ws.Range("D95").Value = "{'andamhie_literal', '(', 'id', 'len', 'korik', 'eme', 'λ', '--', 'chika_literal', '-', '!', '++', 'anda_literal'}"
